$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 10 through 41 (old per-field rows no longer needed)
$ws.Range("A10:A41").EntireRow.Delete() | Out-Null

# Update rows 2-9 with consolidated card tuples
$ws.Range("A2").Value = '(''Blackblade Reforged'', [''{2}'', ''Legendary Artifact — Equipment'', ''Equipped creature gets +1/+1 for each land you control.'', ''Equip legendary creature {3}'', ''Equip {7}''])'
$ws.Range("A3").Value = '(''Gideon Jura'', [''{3}{W}{W}'', ''Legendary Planeswalker — Gideon'', ''+2: During target opponent’s next turn, creatures that player controls attack Gideon Jura if able.'', ''−2: Destroy target tapped creature.'', ''0: Until end of turn, Gideon Jura becomes a 6/6 Human Soldier creature that’s still a planeswalker. Prevent all damage that would be dealt to him this turn.'', ''Loyalty: 6''])'
$ws.Range("A4").Value = '("Martyr''s Bond", [''{4}{W}{W}'', ''Enchantment'', ''Whenever Martyr’s Bond or another nonland permanent you control is put into a graveyard from the battlefield, each opponent sacrifices a permanent that shares a card type with it.''])'
$ws.Range("A5").Value = '(''Path to Exile'', [''{W}'', ''Instant'', ''Exile target creature. Its controller may search their library for a basic land card, put that card onto the battlefield tapped, then shuffle their library.''])'
$ws.Range("A6").Value = '(''Rest in Peace'', [''{1}{W}'', ''Enchantment'', ''When Rest in Peace enters the battlefield, exile all cards from all graveyards.'', ''If a card or token would be put into a graveyard from anywhere, exile it instead.''])'
$ws.Range("A7").Value = '(''Shielded by Faith'', [''{1}{W}{W}'', ''Enchantment — Aura'', ''Enchant creature'', ''Enchanted creature has indestructible.'', ''Whenever a creature enters the battlefield, you may attach Shielded by Faith to that creature.''])'
$ws.Range("A8").Value = '(''True Conviction'', [''{3}{W}{W}{W}'', ''Enchantment'', ''Creatures you control have double strike and lifelink.''])'
$ws.Range("A9").Value = '(''Worship'', [''{3}{W}'', ''Enchantment'', ''If you control a creature, damage that would reduce your life total to less than 1 reduces it to 1 instead.''])'

